# Generate Report for Handback
# Update the Correspond Handoff/Handback datetimes recorded for the
# zh-cn and de-de handback rows (row 2 of each language sheet).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 22:46:23"
$wsZhCn.Range("H2").Value = "2016-03-11 22:46:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 22:46:26"
$wsDeDe.Range("H2").Value = "2016-03-11 22:46:49"
